# This workbook stores a rolling weekly log of price observations.
# A new week of data is being added: the existing block of rows
# (114-146) shifts down by one row (to 115-147, with the old row 147
# observation rolling off the bottom of the window), and a brand new
# observation is written into row 114.
#
# Columns A,B,C,E,F,G,H,I,J,K,T are constant across this block and are
# left untouched; only D (Fecha) and L..S (Calidad .. Precio $/Kg) vary
# per row and need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 114
$lastRow  = 147

# Columns that carry the per-observation data which shifts down by one row.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S")

# Snapshot the current (pre-shift) values for rows firstRow..(lastRow-1)
# before we start overwriting anything. Note: we use .Value2 (not .Value)
# to read/write cell contents throughout this script.
$snapshot = @{}
for ($r = $firstRow; $r -lt $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Push the snapshot down by one row: old row r -> new row r+1.
# Go from the bottom up so we never read a cell after it has been
# overwritten (not strictly required since we used a snapshot, but kept
# for clarity/safety).
for ($r = $lastRow; $r -gt $firstRow; $r--) {
    $src = $snapshot[$r - 1]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $src[$c]
    }
}

# Write the brand new observation into row 114.
$ws.Range("D114").Value2 = 44620
$ws.Range("L114").Value2 = "Primera"
$ws.Range("M114").Value2 = 150
$ws.Range("N114").Value2 = 7000
$ws.Range("O114").Value2 = 7000
$ws.Range("P114").Value2 = 7000
$ws.Range("Q114").Value2 = "$/bandeja 7 kilos"
$ws.Range("R114").Value2 = "Región del Maule"
$ws.Range("S114").Value2 = 1000
